$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("8:10").Delete()

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pspn"
$ws.Range("C2").Value = "Gfra4"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2111796666666667
$ws.Range("H2").Value = 0.633539
$ws.Range("I2").Value = 0.1733132136419605
$ws.Range("J2").Value = 0.1733132136419605
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.600747333333333
$ws.Range("N2").Value = 4.802242
$ws.Range("O2").Value = 0.5378025812999049
$ws.Range("P2").Value = 0.5378025812999049
$ws.Range("Q2").Value = 0.3380452882708889
$ws.Range("R2").Value = 3.042407594438
$ws.Range("S2").Value = 0.09320829367002825
$ws.Range("T2").Value = 0.09320829367002824

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pspn"
$ws.Range("C3").Value = "Gfra4"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2111796666666667
$ws.Range("H3").Value = 0.633539
$ws.Range("I3").Value = 0.1733132136419605
$ws.Range("J3").Value = 0.1733132136419605
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.307232666666666
$ws.Range("N3").Value = 3.921698
$ws.Range("O3").Value = 0.4391905504717742
$ws.Range("P3").Value = 0.4391905504717742
$ws.Range("Q3").Value = 0.2760609588024444
$ws.Range("R3").Value = 2.484548629222
$ws.Range("S3").Value = 0.07611752570344484
$ws.Range("T3").Value = 0.07611752570344484

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pspn"
$ws.Range("C4").Value = "Gfra4"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2111796666666667
$ws.Range("H4").Value = 0.633539
$ws.Range("I4").Value = 0.1733132136419605
$ws.Range("J4").Value = 0.1733132136419605
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.068479
$ws.Range("N4").Value = 0.205437
$ws.Range("O4").Value = 0.02300686822832097
$ws.Range("P4").Value = 0.02300686822832097
$ws.Range("Q4").Value = 0.01446137239366667
$ws.Range("R4").Value = 0.130152351543
$ws.Range("S4").Value = 0.003987394268487426
$ws.Range("T4").Value = 0.003987394268487426

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Pspn"
$ws.Range("C5").Value = "Gfra4"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.007306
$ws.Range("H5").Value = 3.021918
$ws.Range("I5").Value = 0.8266867863580396
$ws.Range("J5").Value = 0.8266867863580396
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.600747333333333
$ws.Range("N5").Value = 4.802242
$ws.Range("O5").Value = 0.5378025812999049
$ws.Range("P5").Value = 0.5378025812999049
$ws.Range("Q5").Value = 1.612442393350667
$ws.Range("R5").Value = 14.511981540156
$ws.Range("S5").Value = 0.4445942876298767
$ws.Range("T5").Value = 0.4445942876298767

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Pspn"
$ws.Range("C6").Value = "Gfra4"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 1.007306
$ws.Range("H6").Value = 3.021918
$ws.Range("I6").Value = 0.8266867863580396
$ws.Range("J6").Value = 0.8266867863580396
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.307232666666666
$ws.Range("N6").Value = 3.921698
$ws.Range("O6").Value = 0.4391905504717742
$ws.Range("P6").Value = 0.4391905504717742
$ws.Range("Q6").Value = 1.316783308529333
$ws.Range("R6").Value = 11.851049776764
$ws.Range("S6").Value = 0.3630730247683294
$ws.Range("T6").Value = 0.3630730247683294

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Pspn"
$ws.Range("C7").Value = "Gfra4"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 1.007306
$ws.Range("H7").Value = 3.021918
$ws.Range("I7").Value = 0.8266867863580396
$ws.Range("J7").Value = 0.8266867863580396
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.068479
$ws.Range("N7").Value = 0.205437
$ws.Range("O7").Value = 0.02300686822832097
$ws.Range("P7").Value = 0.02300686822832097
$ws.Range("Q7").Value = 0.068979307574
$ws.Range("R7").Value = 0.6208137681660001
$ws.Range("S7").Value = 0.01901947395983355
$ws.Range("T7").Value = 0.01901947395983355
